# Gekut met aan elkaar liggende blokken
#
# Appends, after the last paragraph ("0.5 meer brainstormen + testen"):
#   1. An empty (plain, non-list) paragraph tagged en-US
#   2. A plain paragraph tagged nl-NL containing the text
#      "17:20 start 18:45 end (kutten met naast elkaar liggen)"
#      split across three runs, matching the source edit.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Step 1: split a new (empty) paragraph off after "...testen", scoped to
# the last paragraph's own range so we don't clobber the other
# occurrence of "testen" earlier in the document.
# ---------------------------------------------------------------------
$lastPara = $d.Paragraphs.Last
$splitRange1 = $lastPara.Range
$splitRange1.Find.Execute("testen", $true, $false, $false, $false, $false, $true, 1, $false, "testen^p", 2)

# This new trailing paragraph will become the nl-NL text paragraph.
$textPara = $d.Paragraphs.Last
$textPara.Style = "Normal"
$textPara.Range.LanguageID = "nl-NL"

# ---------------------------------------------------------------------
# Step 2: split again at the same anchor to push in a second, empty
# paragraph *before* the nl-NL paragraph created above.
# ---------------------------------------------------------------------
$precedingPara = $d.Paragraphs.Item($d.Paragraphs.Count - 1)
$splitRange2 = $precedingPara.Range
$splitRange2.Find.Execute("testen", $true, $false, $false, $false, $false, $true, 1, $false, "testen^p", 2)

$emptyPara = $d.Paragraphs.Item($d.Paragraphs.Count - 1)
$emptyPara.Style = "Normal"
$emptyPara.Range.LanguageID = "en-US"

# ---------------------------------------------------------------------
# Step 3: fill in the nl-NL paragraph (now the last paragraph) with its
# three runs. Each chunk is typed into its own temporary paragraph (so
# it gets its own run with explicit run-language formatting) and the
# intervening paragraph marks are then deleted to splice the runs back
# together inside a single paragraph, without Word coalescing them into
# one run.
# ---------------------------------------------------------------------
$textPara = $d.Paragraphs.Last

$run1Range = $textPara.Range
$run1Range.Collapse(0)
$run1Range.InsertAfter("17:20 start")
$run1Range.LanguageID = "nl-NL"
$run1Range.InsertParagraphAfter()

$para3 = $d.Paragraphs.Item($d.Paragraphs.Count)
$para3.Range.LanguageID = "nl-NL"
$run2Range = $para3.Range
$run2Range.Collapse(0)
$run2Range.InsertAfter(" 18:45 end (kutten met naast e")
$run2Range.LanguageID = "nl-NL"
$run2Range.InsertParagraphAfter()

$para4 = $d.Paragraphs.Item($d.Paragraphs.Count)
$para4.Range.LanguageID = "nl-NL"
$run3Range = $para4.Range
$run3Range.Collapse(0)
$run3Range.InsertAfter("lkaar liggen)")
$run3Range.LanguageID = "nl-NL"

# Splice the three temporary paragraphs back into one paragraph by
# deleting the paragraph marks between them (re-fetching the paragraph
# boundary fresh each time since indices shift after each delete).
$mergeAnchorIndex = $d.Paragraphs.Count - 2
$mergePara = $d.Paragraphs.Item($mergeAnchorIndex)
$mark1 = $d.Range($mergePara.Range.End - 1, $mergePara.Range.End)
$mark1.Delete()

$mergePara2 = $d.Paragraphs.Item($mergeAnchorIndex)
$mark2 = $d.Range($mergePara2.Range.End - 1, $mergePara2.Range.End)
$mark2.Delete()
